# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column E (estado-de-la-informacion): was a curated dimension (mapped via an
# external mapping file) and becomes a plain measure.
# Column M (municipio-nombre): was a plain measure and becomes a curated
# dimension referencing sdmx-dimension:refArea / URI-Municipio (like
# comarca-nombre already does).
# Column O (tipo-de-presupuesto): was a curated dimension (mapped via an
# external mapping file) and becomes a plain measure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: estado-de-la-informacion -> measure ---
$ws.Range("E2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("E5").Clear()

# --- Column M: municipio-nombre -> dimension (refArea / URI-Municipio) ---
$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("M3").Value = "dim"
$ws.Range("M4").Value = "URI-Municipio"

# --- Column O: tipo-de-presupuesto -> measure ---
$ws.Range("O2").Value = "iaest-measure:tipo-de-presupuesto"
$ws.Range("O3").Value = "medida"
$ws.Range("O4").Value = "xsd:int"
$ws.Range("O5").Clear()
